$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New daily NZ vaccination rows (dates 16-19 Dec 2021), appended after
# the existing last row (302, 15 Dec 2021), bringing the sheet up to
# date with the "21 December" comparison refresh.
$newRows = @(
    @{ Row = 303; Date = 44546; First = 2434; Second = 8294 },
    @{ Row = 304; Date = 44547; First = 2039; Second = 7600 },
    @{ Row = 305; Date = 44548; First = 1728; Second = 6630 },
    @{ Row = 306; Date = 44549; First = 784;  Second = 3360 }
)

foreach ($r in $newRows) {
    $rowNum = $r.Row

    # Copy the formatting of the row above (which already carries the
    # dd/mm/yyyy date style on column A) down onto the new row before
    # writing values, so the new cells pick up the same style index
    # instead of minting a new one.
    $ws.Range("A" + ($rowNum - 1) + ":C" + ($rowNum - 1)).Copy() | Out-Null
    $ws.Range("A" + $rowNum + ":C" + $rowNum).PasteSpecial(-4122) | Out-Null

    $ws.Cells.Item($rowNum, 1).Value = $r.Date
    $ws.Cells.Item($rowNum, 2).Value = $r.First
    $ws.Cells.Item($rowNum, 3).Value = $r.Second
}

$excel.CutCopyMode = 0

# Match the new selection recorded in the saved workbook.
$ws.Range("B304").Select()
